# Rebuild Finance templates with correct industry content
# Replaces "Artificial Intelligence and Machine Learning" project content with
# "Finance - Core Banking System Modernization" content across all tabs.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Instructions & User Guide
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Instructions & User Guide")
$ws1.Range("A1").Value = "Finance - Core Banking System Modernization Comprehensive Budget - User Guide & Instructions"
$ws1.Range("A56").Value = "📋 FINANCE - CORE BANKING SYSTEM MODERNIZATION PROJECT OVERVIEW"

# ---------------------------------------------------------------------------
# Sheet: Budget Summary
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Budget Summary")
$ws2.Range("A1").Value = "Finance - Core Banking System Modernization - Executive Budget Summary"

# ---------------------------------------------------------------------------
# Sheet: Resources
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Resources")
$ws3.Range("A1").Value = "Finance - Core Banking System Modernization - Resources Budget"

# Row 4: Data Scientists -> Banking Systems Architect
$ws3.Range("A4").Value = "Banking Systems Architect"
$ws3.Range("B4").Value = 180
$ws3.Range("C4").Value = 2
$ws3.Range("D4").Value = 14

# Row 5: ML Engineers -> Core Banking Developer
$ws3.Range("A5").Value = "Core Banking Developer"
$ws3.Range("B5").Value = 160
$ws3.Range("C5").Value = 4
$ws3.Range("D5").Value = 14

# Row 6: AI Architects -> Database Administrator
$ws3.Range("A6").Value = "Database Administrator"
$ws3.Range("B6").Value = 150

# Row 7: DevOps Engineers -> Integration Specialist
$ws3.Range("A7").Value = "Integration Specialist"
$ws3.Range("B7").Value = 145
$ws3.Range("C7").Value = 3
$ws3.Range("D7").Value = 14

# Row 8: Project Manager -> QA/Testing Lead
$ws3.Range("A8").Value = "QA/Testing Lead"
$ws3.Range("B8").Value = 130
$ws3.Range("C8").Value = 2
$ws3.Range("D8").Value = 14

# Row 9: Business Analysts -> Business Analyst
$ws3.Range("A9").Value = "Business Analyst"
$ws3.Range("B9").Value = 125
$ws3.Range("C9").Value = 3
$ws3.Range("D9").Value = 14

# Row 10: QA Engineers -> Compliance Officer
$ws3.Range("A10").Value = "Compliance Officer"
$ws3.Range("B10").Value = 140
$ws3.Range("C10").Value = 1
$ws3.Range("D10").Value = 14

# ---------------------------------------------------------------------------
# Sheet: Logistics
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Logistics")
$ws4.Range("A1").Value = "Finance - Core Banking System Modernization - Logistics Budget"
$ws4.Range("B4").Value = 125000
$ws4.Range("B5").Value = 180000
$ws4.Range("B6").Value = 55000
$ws4.Range("B7").Value = 35000
$ws4.Range("B8").Value = 25000

# ---------------------------------------------------------------------------
# Sheet: Technology
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Technology")
$ws5.Range("A1").Value = "Finance - Core Banking System Modernization - Technology Budget"

$ws5.Range("A4").Value = "Core Banking Platform License"
$ws5.Range("B4").Value = 2500000

$ws5.Range("A5").Value = "Cloud Infrastructure (AWS/Azure)"
$ws5.Range("B5").Value = 450000

$ws5.Range("A6").Value = "Database Management System"
$ws5.Range("B6").Value = 280000

$ws5.Range("A7").Value = "Security and Encryption Tools"
$ws5.Range("B7").Value = 180000

$ws5.Range("A8").Value = "API Management Platform"
$ws5.Range("B8").Value = 120000

$ws5.Range("A9").Value = "Testing and QA Tools"
$ws5.Range("B9").Value = 95000

# ---------------------------------------------------------------------------
# Sheet: Training
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Training")
$ws6.Range("A1").Value = "Finance - Core Banking System Modernization - Training Budget"

$ws6.Range("A4").Value = "Core Banking Platform Training"
$ws6.Range("B4").Value = 85000
$ws6.Range("C4").Value = 45

$ws6.Range("A5").Value = "Compliance and Regulatory Training"
$ws6.Range("B5").Value = 45000
$ws6.Range("C5").Value = 20

$ws6.Range("A6").Value = "Technical Skills Development"
$ws6.Range("B6").Value = 65000
$ws6.Range("C6").Value = 30

$ws6.Range("A7").Value = "Change Management Workshops"
$ws6.Range("B7").Value = 35000
$ws6.Range("C7").Value = 45

$ws6.Range("A8").Value = "End-User Training Materials"
$ws6.Range("B8").Value = 25000
$ws6.Range("C8").Value = 200

# ---------------------------------------------------------------------------
# Sheet: Contingency
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("Contingency")
$ws7.Range("A1").Value = "Finance - Core Banking System Modernization - Contingency Budget"
$ws7.Range("D6").Value = "Regulatory changes or compliance requirements"
$ws7.Range("D8").Value = "Staff turnover or skill gaps in banking domain"
$ws7.Range("D9").Value = "Delays or timeline extensions due to testing"

# ---------------------------------------------------------------------------
# Sheet: Timeline
# ---------------------------------------------------------------------------
$ws8 = $wb.Worksheets.Item("Timeline")
$ws8.Range("A1").Value = "Finance - Core Banking System Modernization - Budget Timeline"
